$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 8338.666999999999
$ws.Range("I21").Value = 8338.666999999999
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 8338.666999999999
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -7870.666999999999
$ws.Range("N21").Value = $null

$ws.Range("H23").Value = 8338.666999999999
$ws.Range("I23").Value = 8338.666999999999
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 8338.666999999999
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -8104.666999999999
$ws.Range("N23").Value = $null

$ws.Range("H41").Value = 1316.7778
$ws.Range("I41").Value = 1316.7778
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1316.7778
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -876.7778000000001
$ws.Range("N41").Value = $null

$ws.Range("H42").Value = 612.5
$ws.Range("I42").Value = 612.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 1837.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1607.5
$ws.Range("N42").Value = $null

$ws.Range("H43").Value = 3900
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 5100
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 5100
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -5238

$ws.Range("H53").Value = 1723.5
$ws.Range("I53").Value = 1723.5
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1723.5
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -1086.5
$ws.Range("N53").Value = $null

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = $null

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = $null

$ws.Range("H64").Value = 3946.5
$ws.Range("I64").Value = 3900
$ws.Range("J64").Value = 3993
$ws.Range("K64").Value = 3900
$ws.Range("L64").Value = 3993
$ws.Range("M64").Value = -3652
$ws.Range("N64").Value = -4489

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = $null

$ws.Range("H67").Value = 3946.5
$ws.Range("I67").Value = 3900
$ws.Range("J67").Value = 3993
$ws.Range("K67").Value = 3900
$ws.Range("L67").Value = 3993
$ws.Range("M67").Value = -3042
$ws.Range("N67").Value = -5709

$ws.Range("H69").Value = 13900
$ws.Range("I69").Value = 10500
$ws.Range("J69").Value = 15357.143
$ws.Range("K69").Value = 31500
$ws.Range("L69").Value = 46071.429
$ws.Range("M69").Value = -30626
$ws.Range("N69").Value = -47819.429

$ws.Range("H70").Value = 1832.75
$ws.Range("I70").Value = 1409.4
$ws.Range("J70").Value = 3949.5
$ws.Range("K70").Value = 4228.200000000001
$ws.Range("L70").Value = 11848.5
$ws.Range("M70").Value = -3958.200000000001
$ws.Range("N70").Value = -12388.5

$ws.Range("H72").Value = 13900
$ws.Range("I72").Value = 10500
$ws.Range("J72").Value = 15357.143
$ws.Range("K72").Value = 94500
$ws.Range("L72").Value = 138214.287
$ws.Range("M72").Value = -90132
$ws.Range("N72").Value = -146950.287

$ws.Range("H73").Value = 1832.75
$ws.Range("I73").Value = 1409.4
$ws.Range("J73").Value = 3949.5
$ws.Range("K73").Value = 4228.200000000001
$ws.Range("L73").Value = 11848.5
$ws.Range("M73").Value = -3292.200000000001
$ws.Range("N73").Value = -13720.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = $null

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

$ws.Range("H76").Value = 5265.4443
$ws.Range("I76").Value = 5323.875
$ws.Range("J76").Value = 4798
$ws.Range("K76").Value = 5323.875
$ws.Range("L76").Value = 4798
$ws.Range("M76").Value = -5008.875

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = $null

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

$ws.Range("H79").Value = 5265.4443
$ws.Range("I79").Value = 5323.875
$ws.Range("J79").Value = 4798
$ws.Range("K79").Value = 5323.875
$ws.Range("L79").Value = 4798
$ws.Range("M79").Value = -4231.875

$ws.Range("H80").Value = 2193.9333
$ws.Range("I80").Value = 2354.111
$ws.Range("J80").Value = 1953.6666
$ws.Range("K80").Value = 7062.333
$ws.Range("L80").Value = 5860.9998
$ws.Range("M80").Value = -6064.333
$ws.Range("N80").Value = -7856.9998

$ws.Range("H83").Value = 2193.9333
$ws.Range("I83").Value = 2354.111
$ws.Range("J83").Value = 1953.6666
$ws.Range("K83").Value = 21186.999
$ws.Range("L83").Value = 17582.9994
$ws.Range("M83").Value = -16194.999
$ws.Range("N83").Value = -27566.9994

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = $null

$ws.Range("H100").Value = 4698.7144
$ws.Range("I100").Value = 4698.7144
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4698.7144
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4157.7144

$ws.Range("H116").Value = 5048.4
$ws.Range("I116").Value = 5416.3335
$ws.Range("J116").Value = 4496.5
$ws.Range("K116").Value = 5416.3335
$ws.Range("L116").Value = 4496.5
$ws.Range("M116").Value = -1974.3335

$ws.Range("H127").Value = 2279.625
$ws.Range("I127").Value = 1462.4286
$ws.Range("J127").Value = 8000
$ws.Range("K127").Value = 4387.2858
$ws.Range("L127").Value = 24000
$ws.Range("M127").Value = 572.7142000000003

$ws.Range("H129").Value = 2868.0908
$ws.Range("I129").Value = 954.5
$ws.Range("J129").Value = 3961.5715
$ws.Range("K129").Value = 2863.5
$ws.Range("L129").Value = 11884.7145
$ws.Range("M129").Value = 2136.5

$ws.Range("H131").Value = 5758
$ws.Range("I131").Value = 2577.5
$ws.Range("J131").Value = 9998.666999999999
$ws.Range("K131").Value = 7732.5
$ws.Range("L131").Value = 29996.001
$ws.Range("M131").Value = -2692.5

$ws.Range("H132").Value = 40003748
$ws.Range("I132").Value = 41670424
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 125011272
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -125008742

$ws.Range("H135").Value = 1101.6364
$ws.Range("I135").Value = 1101.6364
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9914.7276
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7379.7276
$ws.Range("N135").Value = $null

$ws.Range("H138").Value = 2185.121
$ws.Range("I138").Value = 2356.2727
$ws.Range("J138").Value = 2099.5454
$ws.Range("K138").Value = 7068.8181
$ws.Range("L138").Value = 6298.6362
$ws.Range("M138").Value = -1928.8181
$ws.Range("N138").Value = -16578.6362

$ws.Range("H141").Value = 11665.667
$ws.Range("I141").Value = 11665.667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 34997.001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -29817.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1654.3636
$ws.Range("I2").Value = 1327.5555
$ws.Range("J2").Value = 3125
$ws.Range("K2").Value = 1327.5555
$ws.Range("L2").Value = 3125
$ws.Range("M2").Value = -1214.5555

$ws.Range("H4").Value = 647.4
$ws.Range("I4").Value = 784.25
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 784.25
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -668.25

$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -713

$ws.Range("H29").Value = 3598.5
$ws.Range("I29").Value = 3598.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3598.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -3290.5

$ws.Range("H32").Value = 6039.2856
$ws.Range("I32").Value = 4546.25
$ws.Range("J32").Value = 14997.5
$ws.Range("K32").Value = 4546.25
$ws.Range("L32").Value = 14997.5
$ws.Range("M32").Value = -4259.25

$ws.Range("H37").Value = 25046
$ws.Range("I37").Value = 2550
$ws.Range("J37").Value = 70038
$ws.Range("K37").Value = 2550
$ws.Range("L37").Value = 70038
$ws.Range("M37").Value = -2277
$ws.Range("N37").Value = -70584

$ws.Range("H45").Value = 5941.778
$ws.Range("I45").Value = 5941.778
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5941.778
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -5564.778

$ws.Range("H61").Value = 1481.95
$ws.Range("I61").Value = 1528.3684
$ws.Range("J61").Value = 600
$ws.Range("K61").Value = 1528.3684
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -1316.3684
$ws.Range("N61").Value = -1024

$ws.Range("H63").Value = 6700
$ws.Range("I63").Value = 1550
$ws.Range("J63").Value = 17000
$ws.Range("K63").Value = 1550
$ws.Range("L63").Value = 17000
$ws.Range("M63").Value = -864
$ws.Range("N63").Value = -18372

$ws.Range("H66").Value = 6700
$ws.Range("I66").Value = 1550
$ws.Range("J66").Value = 17000
$ws.Range("K66").Value = 7750
$ws.Range("L66").Value = 85000
$ws.Range("M66").Value = -4318
$ws.Range("N66").Value = -91864

$ws.Range("H74").Value = 2525.682
$ws.Range("I74").Value = 2021.4286
$ws.Range("J74").Value = 3408.125
$ws.Range("K74").Value = 2021.4286
$ws.Range("L74").Value = 3408.125
$ws.Range("M74").Value = -1147.4286
$ws.Range("N74").Value = -5156.125

$ws.Range("H77").Value = 2525.682
$ws.Range("I77").Value = 2021.4286
$ws.Range("J77").Value = 3408.125
$ws.Range("K77").Value = 10107.143
$ws.Range("L77").Value = 17040.625
$ws.Range("M77").Value = -5739.143
$ws.Range("N77").Value = -25776.625

$ws.Range("H116").Value = 1654.3636
$ws.Range("I116").Value = 1327.5555
$ws.Range("J116").Value = 3125
$ws.Range("K116").Value = 1327.5555
$ws.Range("L116").Value = 3125
$ws.Range("M116").Value = 966.4445000000001

$ws.Range("H122").Value = 769.3333
$ws.Range("I122").Value = 769.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2307.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 142.0001000000002

$ws.Range("H136").Value = 1481.95
$ws.Range("I136").Value = 1528.3684
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 4585.1052
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -2035.1052
$ws.Range("N136").Value = -6900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1654.3636
$ws.Range("I3").Value = 1327.5555
$ws.Range("J3").Value = 3125
$ws.Range("K3").Value = 1327.5555
$ws.Range("L3").Value = 3125
$ws.Range("M3").Value = -1213.5555

$ws.Range("H22").Value = 365.16666
$ws.Range("I22").Value = 365.16666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 365.16666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -192.16666

$ws.Range("H54").Value = 1950
$ws.Range("I54").Value = 1950
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1950
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -1466

$ws.Range("H86").Value = 20032
$ws.Range("I86").Value = 23361.908
$ws.Range("J86").Value = 10874.75
$ws.Range("K86").Value = 23361.908
$ws.Range("L86").Value = 10874.75
$ws.Range("M86").Value = -22238.908

$ws.Range("H89").Value = 20032
$ws.Range("I89").Value = 23361.908
$ws.Range("J89").Value = 10874.75
$ws.Range("K89").Value = 116809.54
$ws.Range("L89").Value = 54373.75
$ws.Range("M89").Value = -111193.54

$ws.Range("H107").Value = 2975
$ws.Range("I107").Value = 2975
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2975
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1055

$ws.Range("H134").Value = 1585.8695
$ws.Range("I134").Value = 1635.2273
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 4905.6819
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = -2370.6819
$ws.Range("N134").Value = -6570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 255.4
$ws.Range("I2").Value = 293
$ws.Range("J2").Value = 199
$ws.Range("K2").Value = 293
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = -180

$ws.Range("H31").Value = 7092.4707
$ws.Range("I31").Value = 18533
$ws.Range("J31").Value = 3572.3076
$ws.Range("K31").Value = 18533
$ws.Range("L31").Value = 3572.3076
$ws.Range("M31").Value = -18238

$ws.Range("H34").Value = 7092.4707
$ws.Range("I34").Value = 18533
$ws.Range("J34").Value = 3572.3076
$ws.Range("K34").Value = 18533
$ws.Range("L34").Value = 3572.3076
$ws.Range("M34").Value = -18331

$ws.Range("H50").Value = 19899

$ws.Range("H51").Value = 38498.5
$ws.Range("I51").Value = 38666.5
$ws.Range("J51").Value = 37994.5
$ws.Range("K51").Value = 38666.5
$ws.Range("L51").Value = 37994.5
$ws.Range("M51").Value = -37930.5
$ws.Range("N51").Value = -39466.5

$ws.Range("H58").Value = 2499.25
$ws.Range("I58").Value = 2499.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2499.25
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2296.25

$ws.Range("H60").Value = 28699.8
$ws.Range("I60").Value = 5750
$ws.Range("J60").Value = 43999.668
$ws.Range("K60").Value = 5750
$ws.Range("L60").Value = 43999.668
$ws.Range("M60").Value = -5239
$ws.Range("N60").Value = -45021.668

$ws.Range("H61").Value = 38498.5
$ws.Range("I61").Value = 38666.5
$ws.Range("J61").Value = 37994.5
$ws.Range("K61").Value = 38666.5
$ws.Range("L61").Value = 37994.5
$ws.Range("M61").Value = -38318.5
$ws.Range("N61").Value = -38690.5

$ws.Range("H62").Value = 3225
$ws.Range("I62").Value = 3225
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3225
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2601

$ws.Range("H65").Value = 3225
$ws.Range("I65").Value = 3225
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16125
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13005

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = $null

$ws.Range("H136").Value = 2499.25
$ws.Range("I136").Value = 2499.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7497.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4947.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1598.4
$ws.Range("I11").Value = 1664.6666
$ws.Range("J11").Value = 1499
$ws.Range("K11").Value = 4993.9998
$ws.Range("L11").Value = 4497
$ws.Range("M11").Value = -4853.9998
$ws.Range("N11").Value = -4777

$ws.Range("H12").Value = 90.666664
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 88.8
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 266.4
$ws.Range("M12").Value = -127
$ws.Range("N12").Value = -612.4

$ws.Range("H56").Value = 9290.286
$ws.Range("I56").Value = 9290.286
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 9290.286
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8760.286

$ws.Range("H68").Value = 727.7778
$ws.Range("I68").Value = 739.2
$ws.Range("J68").Value = 713.5
$ws.Range("K68").Value = 2217.6
$ws.Range("L68").Value = 2140.5
$ws.Range("M68").Value = -1406.6
$ws.Range("N68").Value = -3762.5

$ws.Range("H71").Value = 727.7778
$ws.Range("I71").Value = 739.2
$ws.Range("J71").Value = 713.5
$ws.Range("K71").Value = 6652.8
$ws.Range("L71").Value = 6421.5
$ws.Range("M71").Value = -2596.8
$ws.Range("N71").Value = -14533.5

$ws.Range("H94").Value = 11791.5
$ws.Range("I94").Value = 750
$ws.Range("J94").Value = 13999.8
$ws.Range("K94").Value = 2250
$ws.Range("L94").Value = 41999.39999999999
$ws.Range("M94").Value = -1574
$ws.Range("N94").Value = -43351.39999999999

$ws.Range("H129").Value = 1287
$ws.Range("I129").Value = 1081.6666
$ws.Range("J129").Value = 1441
$ws.Range("K129").Value = 3244.9998
$ws.Range("L129").Value = 4323
$ws.Range("M129").Value = 1755.0002
$ws.Range("N129").Value = -14323

$ws.Range("H131").Value = 1829.1082
$ws.Range("I131").Value = 4000
$ws.Range("J131").Value = 1768.8055
$ws.Range("K131").Value = 12000
$ws.Range("L131").Value = 5306.416499999999
$ws.Range("M131").Value = -6960
$ws.Range("N131").Value = -15386.4165

$ws.Range("H139").Value = 7338.423
$ws.Range("I139").Value = 5133.1665
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 15399.4995
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = -10259.4995

$ws.Range("H140").Value = 1113398.5
$ws.Range("I140").Value = 1113398.5
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3340195.5
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -3335015.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2454.2
$ws.Range("I102").Value = 1893.5555
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 1893.5555
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -271.5554999999999

$ws.Range("H107").Value = 2569.5715
$ws.Range("I107").Value = 799
$ws.Range("J107").Value = 6996
$ws.Range("K107").Value = 799
$ws.Range("L107").Value = 6996
$ws.Range("M107").Value = 1121
$ws.Range("N107").Value = -10836

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3894.5
$ws.Range("I7").Value = 3894.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3894.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3782.5
$ws.Range("N7").Value = $null

$ws.Range("H22").Value = 41667890
$ws.Range("I22").Value = 986.6667
$ws.Range("J22").Value = 55556856
$ws.Range("K22").Value = 986.6667
$ws.Range("L22").Value = 55556856
$ws.Range("M22").Value = -691.6667
$ws.Range("N22").Value = -55557446

$ws.Range("H27").Value = 41667890
$ws.Range("I27").Value = 986.6667
$ws.Range("J27").Value = 55556856
$ws.Range("K27").Value = 986.6667
$ws.Range("L27").Value = 55556856
$ws.Range("M27").Value = -879.6667
$ws.Range("N27").Value = -55557070

$ws.Range("H46").Value = 3677.7334
$ws.Range("I46").Value = 1104.75
$ws.Range("J46").Value = 4613.364
$ws.Range("K46").Value = 1104.75
$ws.Range("L46").Value = 4613.364
$ws.Range("M46").Value = -916.75
$ws.Range("N46").Value = -4989.364

$ws.Range("H68").Value = 2438.25
$ws.Range("I68").Value = 1375
$ws.Range("J68").Value = 3501.5
$ws.Range("K68").Value = 1375
$ws.Range("L68").Value = 3501.5
$ws.Range("M68").Value = -626
$ws.Range("N68").Value = -4999.5

$ws.Range("H71").Value = 2438.25
$ws.Range("I71").Value = 1375
$ws.Range("J71").Value = 3501.5
$ws.Range("K71").Value = 6875
$ws.Range("L71").Value = 17507.5
$ws.Range("M71").Value = -3131
$ws.Range("N71").Value = -24995.5

$ws.Range("H126").Value = 3894.5
$ws.Range("I126").Value = 3894.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11683.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9213.5
$ws.Range("N126").Value = $null

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 57500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 57500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 57500
$ws.Range("N80").Value = -59496

$ws.Range("H82").Value = 44966.668
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 62450
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 62450
$ws.Range("M82").Value = -9617
$ws.Range("N82").Value = -63216

$ws.Range("H83").Value = 57500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 57500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 172500
$ws.Range("N83").Value = -182484

$ws.Range("H85").Value = 44966.668
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 62450
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 62450
$ws.Range("M85").Value = -8674
$ws.Range("N85").Value = -65102

$ws.Range("H88").Value = 16000
$ws.Range("I88").Value = 16000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 16000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -15594

$ws.Range("H91").Value = 16000
$ws.Range("I91").Value = 16000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 16000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -14596

$ws.Range("H96").Value = 3466.3333
$ws.Range("I96").Value = 4599.5
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 4599.5
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -3226.5
$ws.Range("N96").Value = -3946

$ws.Range("H105").Value = 12790
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12790
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 12790
$ws.Range("M105").Value = $null
$ws.Range("N105").Value = -19778

$ws.Range("H122").Value = 1201.125
$ws.Range("I122").Value = 1287
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3861
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -1411

$ws.Range("H123").Value = 59993.332
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 59993.332
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 59993.332
$ws.Range("N123").Value = -69793.33199999999

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = $null

$ws.Range("H132").Value = 1000000000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1000000000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3000000000
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -3000005060

$ws.Range("H138").Value = 73107.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 73107.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 73107.5
$ws.Range("N138").Value = -83387.5

$ws.Range("H141").Value = 72500.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 72500.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 72500.5
$ws.Range("N141").Value = -82860.5
